$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D24").Value = "[논문 요약] Few-Shot Parameter-Efficient Fine-Tuning is Better and Cheaper than In-Context Learning"
$ws.Range("E24").Value = "https://blog.naver.com/hist0134/222752748525"

$ws.Range("D42").Value = "python pip mirror 서버 설정"
$ws.Range("E42").Value = "https://kjk92.tistory.com/83"

$ws.Range("D51").Value = "[MS Excel 2010] 중복된 데이터에서 고유값을 남겨 놓고 각 그룹의 데이터 개수, 평균, 최대값, 최소값, 중앙값, 표준편차 구하기"
$ws.Range("E51").Value = "https://bskyvision.com/1283"
